$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row.
# Update it from 2023-10-07 (45206) to 2023-10-09 (45208) for all data rows.
$newDateSerial = 45208
$firstRow = 2
$lastRow = 469

$ws.Range($ws.Cells.Item($firstRow, 3), $ws.Cells.Item($lastRow, 3)).Value = $newDateSerial
